$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20
$data[0,0] = "ECs"
$data[0,1] = "Efna1"
$data[0,2] = "Epha3"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 31.824752
$data[0,7] = 95.47425600000001
$data[0,8] = 0.886907633630525
$data[0,9] = 0.886907633630525
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.009849666666666666
$data[0,13] = 0.029549
$data[0,14] = 0.0002013876315934659
$data[0,15] = 0.0002013876315934659
$data[0,16] = 0.3134631989493333
$data[0,17] = 2.821168790544
$data[0,18] = 0.0001786122277790168
$data[0,19] = 0.0001786122277790168
$data[1,0] = "ECs"
$data[1,1] = "Efna1"
$data[1,2] = "Epha3"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 31.824752
$data[1,7] = 95.47425600000001
$data[1,8] = 0.886907633630525
$data[1,9] = 0.886907633630525
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 46.25093466666667
$data[1,13] = 138.752804
$data[1,14] = 0.9456529349389956
$data[1,15] = 0.9456529349389956
$data[1,16] = 1471.924525534869
$data[1,17] = 13247.32072981383
$data[1,18] = 0.8387068067625054
$data[1,19] = 0.8387068067625054
$data[2,0] = "ECs"
$data[2,1] = "Efna1"
$data[2,2] = "Epha3"
$data[2,3] = "MuSCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 31.824752
$data[2,7] = 95.47425600000001
$data[2,8] = 0.886907633630525
$data[2,9] = 0.886907633630525
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 2.607896333333333
$data[2,13] = 7.823689
$data[2,14] = 0.05332140505715427
$data[2,15] = 0.05332140505715428
$data[2,16] = 82.99565405004267
$data[2,17] = 746.9608864503841
$data[2,18] = 0.04729116118109541
$data[2,19] = 0.04729116118109541
$data[3,0] = "ECs"
$data[3,1] = "Efna1"
$data[3,2] = "Epha3"
$data[3,3] = "Resolving-Mac"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 31.824752
$data[3,7] = 95.47425600000001
$data[3,8] = 0.886907633630525
$data[3,9] = 0.886907633630525
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.04031433333333333
$data[3,13] = 0.120943
$data[3,14] = 0.0008242723722565416
$data[3,15] = 0.0008242723722565415
$data[3,16] = 1.282993660378667
$data[3,17] = 11.546942943408
$data[3,18] = 0.0007310534591450686
$data[3,19] = 0.0007310534591450685
$data[4,0] = "FAPs"
$data[4,1] = "Efna1"
$data[4,2] = "Epha3"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 2.333117
$data[4,7] = 6.999351000000001
$data[4,8] = 0.06502043684278042
$data[4,9] = 0.06502043684278042
$data[4,10] = 1
$data[4,11] = 0.3333333333333333
$data[4,12] = 0.009849666666666666
$data[4,13] = 0.029549
$data[4,14] = 0.0002013876315934659
$data[4,15] = 0.0002013876315934659
$data[4,16] = 0.02298042474433333
$data[4,17] = 0.206823822699
$data[4,18] = 0.00001309431178094008
$data[4,19] = 0.00001309431178094008
$data[5,0] = "FAPs"
$data[5,1] = "Efna1"
$data[5,2] = "Epha3"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 2.333117
$data[5,7] = 6.999351000000001
$data[5,8] = 0.06502043684278042
$data[5,9] = 0.06502043684278042
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 46.25093466666667
$data[5,13] = 138.752804
$data[5,14] = 0.9456529349389956
$data[5,15] = 0.9456529349389956
$data[5,16] = 107.9088419366893
$data[5,17] = 971.1795774302041
$data[5,18] = 0.0614867669313909
$data[5,19] = 0.0614867669313909
$data[6,0] = "FAPs"
$data[6,1] = "Efna1"
$data[6,2] = "Epha3"
$data[6,3] = "MuSCs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 2.333117
$data[6,7] = 6.999351000000001
$data[6,8] = 0.06502043684278042
$data[6,9] = 0.06502043684278042
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 2.607896333333333
$data[6,13] = 7.823689
$data[6,14] = 0.05332140505715427
$data[6,15] = 0.05332140505715428
$data[6,16] = 6.084527269537666
$data[6,17] = 54.760745425839
$data[6,18] = 0.003466981049887011
$data[6,19] = 0.003466981049887012
$data[7,0] = "FAPs"
$data[7,1] = "Efna1"
$data[7,2] = "Epha3"
$data[7,3] = "Resolving-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 2.333117
$data[7,7] = 6.999351000000001
$data[7,8] = 0.06502043684278042
$data[7,9] = 0.06502043684278042
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.04031433333333333
$data[7,13] = 0.120943
$data[7,14] = 0.0008242723722565416
$data[7,15] = 0.0008242723722565415
$data[7,16] = 0.09405805644366667
$data[7,17] = 0.8465225079930001
$data[7,18] = 0.00005359454972155525
$data[7,19] = 0.00005359454972155525
$data[8,0] = "MuSCs"
$data[8,1] = "Efna1"
$data[8,2] = "Epha3"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 1.696588
$data[8,7] = 5.089764000000001
$data[8,8] = 0.04728133775640876
$data[8,9] = 0.04728133775640876
$data[8,10] = 1
$data[8,11] = 0.3333333333333333
$data[8,12] = 0.009849666666666666
$data[8,13] = 0.029549
$data[8,14] = 0.0002013876315934659
$data[8,15] = 0.0002013876315934659
$data[8,16] = 0.01671082627066667
$data[8,17] = 0.150397436436
$data[8,18] = 0.000009521876629333878
$data[8,19] = 0.000009521876629333878
$data[9,0] = "MuSCs"
$data[9,1] = "Efna1"
$data[9,2] = "Epha3"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 1.696588
$data[9,7] = 5.089764000000001
$data[9,8] = 0.04728133775640876
$data[9,9] = 0.04728133775640876
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 46.25093466666667
$data[9,13] = 138.752804
$data[9,14] = 0.9456529349389956
$data[9,15] = 0.9456529349389956
$data[9,16] = 78.46878074425068
$data[9,17] = 706.2190266982561
$data[9,18] = 0.04471173581718989
$data[9,19] = 0.04471173581718989
$data[10,0] = "MuSCs"
$data[10,1] = "Efna1"
$data[10,2] = "Epha3"
$data[10,3] = "MuSCs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 1.696588
$data[10,7] = 5.089764000000001
$data[10,8] = 0.04728133775640876
$data[10,9] = 0.04728133775640876
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 2.607896333333333
$data[10,13] = 7.823689
$data[10,14] = 0.05332140505715427
$data[10,15] = 0.05332140505715428
$data[10,16] = 4.424525624377334
$data[10,17] = 39.820730619396
$data[10,18] = 0.002521107362153593
$data[10,19] = 0.002521107362153593
$data[11,0] = "MuSCs"
$data[11,1] = "Efna1"
$data[11,2] = "Epha3"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 1.696588
$data[11,7] = 5.089764000000001
$data[11,8] = 0.04728133775640876
$data[11,9] = 0.04728133775640876
$data[11,10] = 1
$data[11,11] = 0.3333333333333333
$data[11,12] = 0.04031433333333333
$data[11,13] = 0.120943
$data[11,14] = 0.0008242723722565416
$data[11,15] = 0.0008242723722565415
$data[11,16] = 0.06839681416133335
$data[11,17] = 0.615571327452
$data[11,18] = 0.00003897270043593784
$data[11,19] = 0.00003897270043593783
$data[12,0] = "Resolving-Mac"
$data[12,1] = "Efna1"
$data[12,2] = "Epha3"
$data[12,3] = "ECs"
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.02836866666666667
$data[12,7] = 0.085106
$data[12,8] = 0.0007905917702857979
$data[12,9] = 0.0007905917702857978
$data[12,10] = 1
$data[12,11] = 0.3333333333333333
$data[12,12] = 0.009849666666666666
$data[12,13] = 0.029549
$data[12,14] = 0.0002013876315934659
$data[12,15] = 0.0002013876315934659
$data[12,16] = 0.0002794219104444444
$data[12,17] = 0.002514797194
$data[12,18] = 0.0000001592154041751423
$data[12,19] = 0.0000001592154041751423
$data[13,0] = "Resolving-Mac"
$data[13,1] = "Efna1"
$data[13,2] = "Epha3"
$data[13,3] = "FAPs"
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.02836866666666667
$data[13,7] = 0.085106
$data[13,8] = 0.0007905917702857979
$data[13,9] = 0.0007905917702857978
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 46.25093466666667
$data[13,13] = 138.752804
$data[13,14] = 0.9456529349389956
$data[13,15] = 0.9456529349389956
$data[13,16] = 1.312077348580444
$data[13,17] = 11.808696137224
$data[13,18] = 0.000747625427909381
$data[13,19] = 0.0007476254279093809
$data[14,0] = "Resolving-Mac"
$data[14,1] = "Efna1"
$data[14,2] = "Epha3"
$data[14,3] = "MuSCs"
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.02836866666666667
$data[14,7] = 0.085106
$data[14,8] = 0.0007905917702857979
$data[14,9] = 0.0007905917702857978
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 2.607896333333333
$data[14,13] = 7.823689
$data[14,14] = 0.05332140505715427
$data[14,15] = 0.05332140505715428
$data[14,16] = 0.07398254178155556
$data[14,17] = 0.6658428760340001
$data[14,18] = 0.00004215546401826169
$data[14,19] = 0.00004215546401826169
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Efna1"
$data[15,2] = "Epha3"
$data[15,3] = "Resolving-Mac"
$data[15,4] = 1
$data[15,5] = 0.3333333333333333
$data[15,6] = 0.02836866666666667
$data[15,7] = 0.085106
$data[15,8] = 0.0007905917702857979
$data[15,9] = 0.0007905917702857978
$data[15,10] = 1
$data[15,11] = 0.3333333333333333
$data[15,12] = 0.04031433333333333
$data[15,13] = 0.120943
$data[15,14] = 0.0008242723722565416
$data[15,15] = 0.0008242723722565415
$data[15,16] = 0.001143663884222222
$data[15,17] = 0.010292974958
$data[15,18] = 0.0000006516629539799734
$data[15,19] = 0.0000006516629539799733

$ws.Range("A2:T17").Value = $data

Write-Output "done"